# Apply "formatted ToDo list and added env.yml file" changes to config.xlsx
# (the env.yml file itself lives outside this workbook; only the
# spreadsheet-visible edits -- two new simulation columns (G, H) plus
# renamed/updated group + coordinate data -- are applied here)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: missile_name -- add two new simulations (columns G & H) ---
$ws.Range("G2").Value = "missile4"
$ws.Range("H2").Value = "missile5"

# --- Row 3: group_name -- rename existing groups and extend to G & H ---
$ws.Range("D3").Value = "test_group1"
$ws.Range("E3").Value = "test_group1"
$ws.Range("F3").Value = "test_group2"
$ws.Range("G3").Value = "test_group2"
$ws.Range("H3").Value = "test_group2"

# --- Row 4: launch_date -- extend to G & H (serial date 12/6/2020) ---
$ws.Range("G4").Value2 = 44042
$ws.Range("H4").Value2 = 44042

# --- Row 5: launch_time_UTC -- per-simulation launch times ---
$ws.Range("E5").Value2 = 0.16701388888888891
$ws.Range("F5").Value2 = 0.16718750000000002
$ws.Range("G5").Value2 = 0.16770833333333335
$ws.Range("H5").Value2 = 0.16805555555555554

# --- Row 6: LP_lat_deg ---
$ws.Range("D6").Value2 = 39.516824999999997
$ws.Range("E6").Value2 = 39.485511000000002
$ws.Range("F6").Value2 = 39.504075999999998
$ws.Range("G6").Value2 = 39.542740999999999
$ws.Range("H6").Value2 = 39.521768999999999

# --- Row 7: LP_lon_deg ---
$ws.Range("D7").Value2 = -104.95567
$ws.Range("E7").Value2 = -104.884624
$ws.Range("F7").Value2 = -104.90361900000001
$ws.Range("G7").Value2 = -104.97689699999999
$ws.Range("H7").Value2 = -105.01303799999999

# --- Row 8: AP_lat_deg ---
$ws.Range("D8").Value2 = 40.862397000000001
$ws.Range("E8").Value2 = 40.862397000000001
$ws.Range("F8").Value2 = 40.862397000000001
$ws.Range("G8").Value2 = 40.862397000000001
$ws.Range("H8").Value2 = 40.862397000000001
$ws.Range("G8").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("H8").NumberFormat = $ws.Range("D8").NumberFormat

# --- Row 9: AP_lon_deg ---
$ws.Range("D9").Value2 = -105.025902
$ws.Range("E9").Value2 = -105.025902
$ws.Range("F9").Value2 = -105.025902
$ws.Range("G9").Value2 = -105.025902
$ws.Range("H9").Value2 = -105.025902
$ws.Range("G9").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("H9").NumberFormat = $ws.Range("D9").NumberFormat

# --- Row 10: horizontal_velocity_km_sec ---
$ws.Range("G10").Value2 = 1
$ws.Range("H10").Value2 = 1

# --- Row 11: timestep_sec ---
$ws.Range("G11").Value2 = 1
$ws.Range("H11").Value2 = 1

# --- Row 12: sim_start_time_buffer_sec ---
$ws.Range("G12").Value2 = 10
$ws.Range("H12").Value2 = 10

# --- Row 13: sim_end_time_buffer_sec ---
$ws.Range("G13").Value2 = 10
$ws.Range("H13").Value2 = 10

# --- Row 14: collada_model_dir ---
$ws.Range("G14").Value = "../Blender"
$ws.Range("H14").Value = "../Blender"

# --- Row 15: collada_model_file ---
$ws.Range("G15").Value = "test_missile.dae"
$ws.Range("H15").Value = "test_missile.dae"

# --- Row 16: collada_model_scale ---
$ws.Range("G16").Value2 = 500
$ws.Range("H16").Value2 = 500

# --- Selection moved to H7 ---
$ws.Range("H7").Select() | Out-Null
